# Adding the changes we made on may 9th
#
# Semantics (derived from the target diff):
#   - 18 brand-new data rows are inserted directly under the header row,
#     pushing the existing data down.
#   - The 8 oldest data rows (which fall off the bottom once the new rows
#     are inserted) are removed, leaving the first 12 of the original 20
#     data rows intact at the bottom of the sheet.
#   - Net effect: sheet grows from 20 data rows (A1:C21) to 30 data rows
#     (A1:C31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to insert right after the header (row 1).
$newRows = @(
    @(-0.0005742134153842587, 0.0210748501121997, 0.04580267548561091),
    @(-0.0293826170265674, 0.0331699833273887, 0.01545489076524965),
    @(-0.02163684628903856, 0.01164309203624722, -0.01903456177562478),
    @(-0.00761748962104312, -0.02296242564916602, -0.00737925238907336),
    @(-0.03987117595970628, -0.01981036089360703, 0.01966986127197736),
    @(-0.04132503550499657, -0.01673159938305617, 0.0125227374956011),
    @(-0.01340849157422786, -0.03075706511735912, -0.009834930114448132),
    @(-0.03994447708129877, -0.01221119597554203, -0.04150218397378917),
    @(-0.03637702405452719, 0.01701259657740584, -0.03676186949014657),
    @(-0.03740938737988467, 0.01381166309118263, -0.03841120541095729),
    @(-0.04191757388412946, 0.01410487815737716, -0.03918089691549535),
    @(0.01087340153753754, 0.01499674115329971, -0.01282817013561728),
    @(0.02458121769130225, 0.01081842321902504, -0.04899139240384098),
    @(0.01802052438259117, 0.00392786357551804, -0.03619987547397605),
    @(0.01199739351868624, -0.0002138027921320222, 0.03830125063657756),
    @(0.01877188928425304, -0.01797165483236304, 0.03120299618691198),
    @(-0.0108428578823804, -0.0058032199740409, -0.0042760567739605),
    @(0.01026864476501944, 0.01618792921304707, 0.006133087240159481)
)

$insertCount = $newRows.Count

# Insert that many blank rows directly below the header row (row 1); this
# shoves the existing data (old rows 2-21) down to rows (2+insertCount)..(21+insertCount).
$insertRange = $ws.Range("A2:A" + ($insertCount + 1))
$insertRange.EntireRow.Insert()

# The inserted rows inherit the header's (bold) formatting from the row
# above; strip that back off so the new rows look like ordinary data rows.
$ws.Range("A2:C" + ($insertCount + 1)).ClearFormats()

# Fill the freshly inserted rows with the new data.
for ($i = 0; $i -lt $insertCount; $i++) {
    $row = $newRows[$i]
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# The last 8 rows of the (now shifted-down) original data fall away: they
# used to be old rows 14-21, which are now at rows (14+insertCount)..(21+insertCount).
$firstDropRow = 14 + $insertCount
$lastDropRow = 21 + $insertCount
$ws.Range("A" + $firstDropRow + ":A" + $lastDropRow).EntireRow.Delete() | Out-Null

Write-Host "Rows inserted: $insertCount; dropped rows $firstDropRow-$lastDropRow"
